# Insert a new row (ID 3532 / "금호베스트빌") above the existing row 33
# ("서울숲대림", ID 3530) and shift every following row down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Insert()
$ws.Range("A33").Value = 3532
$ws.Range("B33").Value = "금호베스트빌"

# The original sheet applied a (now-unused) number-format style to every
# cell in column A (s="1"). The authored edit drops that style entirely,
# so clear formatting on the whole ID column (including the freshly
# inserted row, which otherwise inherits formatting from the row above).
$ws.Columns.Item(1).ClearFormats()

# Match the author's final selection/viewport.
$ws.Range("B34").Select() | Out-Null
